# Apply the weekly price-record rotation described in the commit
# "Fruta / hortaliza, semanal": each listed row's observation data
# (Fecha, Volumen, Precio minimo/maximo/promedio, Origen, Precio $/Kg)
# is replaced by another week's record, per the verified target diff.
# Identity columns (Mercado, Region, Codreg, Tipo, Producto, Categoria,
# Variedad, Calidad, Unidad, Kg/unidad) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D3").Value = 44582
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6500
$ws.Range("P3").Value = 6233
$ws.Range("S3").Value = 3116

$ws.Range("D5").Value = 44588
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 6500
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 6750
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 3375

$ws.Range("D6").Value = 44586
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 7000
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 3500

$ws.Range("D7").Value = 44585
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 6500
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 6750
$ws.Range("S7").Value = 3375

$ws.Range("D8").Value = 44589
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 6000
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 3000

$ws.Range("D10").Value = 44214
$ws.Range("M10").Value = 48
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("S10").Value = 3000

$ws.Range("D11").Value = 44209
$ws.Range("M11").Value = 58
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 6000
$ws.Range("S11").Value = 3000

$ws.Range("D12").Value = 44587
$ws.Range("M12").Value = 165
$ws.Range("N12").Value = 6500
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 6742
$ws.Range("S12").Value = 3371

$ws.Range("D13").Value = 44606
$ws.Range("M13").Value = 45
$ws.Range("R13").Value = "Provincia de Linares"

$ws.Range("D14").Value = 44614
$ws.Range("M14").Value = 45
$ws.Range("N14").Value = 6000
$ws.Range("O14").Value = 6000
$ws.Range("P14").Value = 6000
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value = 3000

$ws.Range("D15").Value = 44627
$ws.Range("M15").Value = 45
$ws.Range("R15").Value = "Provincia de Linares"

